$d = $word.ActiveDocument

# The "_GoBack" bookmark sits in an otherwise-empty ListParagraph
# (ind left=2160) that immediately precedes the "Evaluate potential
# Solution" / "Does each solution meet the goals?" / "Will each solution
# work for ALL cases?" bullets belonging to the sock-drawer problem
# (numId 15). Locate that anchor paragraph by scanning the document's
# Paragraphs collection for the one whose range contains the bookmark.
#
# NOTE: reading .Text/.Paragraphs off a *manually constructed* zero-length
# Document.Range(x, x) is unreliable in this host, so all reads below go
# through Paragraph objects obtained from the Paragraphs collection
# (whose own .Range is safe to read). Document.Range(x, x) is only used
# as a write-target for InsertXML, which works correctly with raw
# integer offsets.
$bm = $d.Bookmarks.Item("_GoBack")
$bmStart = $bm.Start

$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $bmStart -and $bmStart -lt $p.Range.End) {
        $anchorIdx = $i
        break
    }
}

$pGoals = $d.Paragraphs.Item($anchorIdx + 2)
$pWorks = $d.Paragraphs.Item($anchorIdx + 3)

# Reword the two sub-bullets.
$pGoals.Range.Text = "The solution meets goals"
$pWorks.Range.Text = "The solution works for all cases."

# Insert a new blank ListParagraph (ind left=1440, no numbering) right
# after "The solution works for all cases." -- this paragraph will become
# the new home of the "_GoBack" bookmark.
$insertPos = $pWorks.Range.End
$insertRange = $d.Range($insertPos, $insertPos)
$blankParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($blankParaXml)

# Move the "_GoBack" bookmark from its old (now plain) paragraph to the
# newly-created blank paragraph.
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

$newBlankPara = $d.Paragraphs.Item($anchorIdx + 4)
$d.Bookmarks.Add("_GoBack", $newBlankPara.Range)
